$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "-"
$ws.Range("B3").Value = "MCT-2A-M.T.R.M."
$ws.Range("C3").Value = "[-, -, 'MCT-3A-Elementos de máquinas', -]"
$ws.Range("B4").Value = "MCT-2A-M.T.R.M."
$ws.Range("C4").Value = "[-, -, 'MCT-3A-Elementos de máquinas', -]"
$ws.Range("E4").Value = "-"
$ws.Range("C6").Value = "[-, -, 'MCT-3A-Elementos de máquinas', -]"
$ws.Range("E6").Value = "-"
$ws.Range("C7").Value = "[-, -, 'MCT-3A-Elementos de máquinas', -]"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "-"
$ws.Range("D8").Value = "-"
